# "add select2 to readme"
# The "External Sources" paragraph lists third-party tools used in the
# project. The commit adds the Select2 widget (and the jQuery it depends
# on) to that list, replacing the previous lone mention of "Ajax".
#
# The word "Ajax" lives in its own w:r (run), sandwiched between a run
# ending in "...does not exist; " and a run that is just the trailing
# ".". We only want to retarget the text of that middle run, leaving the
# neighboring runs completely untouched/unmerged.

$d = $word.ActiveDocument
$replacement = "Jquery; Select2 widget"

$r = $d.Content
$found = $r.Find.Execute("Ajax", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertStart = $r.Start

    # Replace the text of just this run.
    $r.Text = $replacement

    # Re-acquire the range covering the text we just inserted and nudge
    # its direct character formatting (on, then back off) so the engine
    # keeps it as its own run instead of silently coalescing it with the
    # neighboring runs that happen to share the same (empty) formatting.
    $new = $d.Range($insertStart, $insertStart + $replacement.Length)
    $new.Font.Bold = 1
    $new.Font.Bold = 0
}
